$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh latest crypto prices / 1h volume change percentages
# (rows 47-48: Aave/Frax swapped places in the ranking)

# Row 2
$ws.Range("D2").Value = '25.919.79'
$ws.Range("E2").Value = '  +0.24%  '

# Row 3
$ws.Range("D3").Value = '1.648.45'
$ws.Range("E3").Value = '  +0.77%  '

# Row 4
$ws.Range("E4").Value = '  +0.66%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.74'
$ws.Range("E5").Value = '  +0.01%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5091'
$ws.Range("E6").Value = '  +1.25%  '

# Row 7
$ws.Range("E7").Value = '  +0.48%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2579'
$ws.Range("E8").Value = '  +0.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06423'
$ws.Range("E9").Value = '  +0.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.72'
$ws.Range("E10").Value = '  +0.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07780'
$ws.Range("E11").Value = '  +1.15%  '

# Row 12
$ws.Range("D12").Value = '1.678.85'
$ws.Range("E12").Value = '  +2.97%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.317'
$ws.Range("E13").Value = '  +1.71%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5475'
$ws.Range("E14").Value = '  +0.40%  '

# Row 15
$ws.Range("D15").Value = '0.0₅7910'
$ws.Range("E15").Value = '  -0.44%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.10'
$ws.Range("E16").Value = '  +2.47%  '

# Row 17
$ws.Range("D17").Value = '26.012.01'
$ws.Range("E17").Value = '  +0.60%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.007'
$ws.Range("E18").Value = '  +0.50%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '197.64'
$ws.Range("E19").Value = '  -2.58%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.434'
$ws.Range("E20").Value = '  +2.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.04'
$ws.Range("E21").Value = '  +0.89%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.076'
$ws.Range("E22").Value = '  +1.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.009'
$ws.Range("E23").Value = '  +0.61%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.861'
$ws.Range("E24").Value = '  -3.93%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.14'
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1148'
$ws.Range("E26").Value = '  -0.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.899'
$ws.Range("E27").Value = '  +2.76%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.77'
$ws.Range("E28").Value = '  +0.52%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.243'
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05022'
$ws.Range("E30").Value = '  +0.02%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.279'
$ws.Range("E31").Value = '  +0.14%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.204'
$ws.Range("E32").Value = '  +0.54%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.546'
$ws.Range("E33").Value = '  +0.71%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.370'
$ws.Range("E34").Value = '  +0.45%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8960'
$ws.Range("E35").Value = '  +0.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.592'
$ws.Range("E36").Value = '  -0.59%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5563'
$ws.Range("E37").Value = '  -0.92%  '

# Row 38
$ws.Range("D38").Value = '1.133.81'
$ws.Range("E38").Value = '  -3.22%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01568'
$ws.Range("E39").Value = '  +0.57%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.008'
$ws.Range("E40").Value = '  +0.68%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.670'
$ws.Range("E41").Value = '  +0.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8171'
$ws.Range("E42").Value = '  +1.24%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.84'
$ws.Range("E43").Value = '  +0.24%  '

# Row 44
$ws.Range("D44").Value = '0.0₈125'
$ws.Range("E44").Value = '  +7.99%  '

# Row 45
$ws.Range("D45").Value = '1.785.14'
$ws.Range("E45").Value = '  +0.78%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4540'
$ws.Range("E46").Value = '  +0.64%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.32'
$ws.Range("E47").Value = '  +1.00%  '

# Row 48
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  +0.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05095'
$ws.Range("E49").Value = '  +0.90%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.007'
$ws.Range("E50").Value = '  +0.39%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09556'
$ws.Range("E51").Value = '  +2.89%  '
